$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D8","D10","D11","D13","D15","D18","D19","D22","D25","D27","D31","D36","D39","D40","D44","D47","D48","D49")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range('D2').Value = '34.125.68'
$ws.Range('D3').Value = '1.791.04'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '227.13'
$ws.Range('E5').Value = '  +1.16%  '
$ws.Range('E6').Value = '  -0.68%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '32.43'
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('E9').Value = '  +4.33%  '
$ws.Range('D10').Value = '0.0689'
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('D11').Value = '0.0941'
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').Value = '2.048.93'
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').Value = '11.49'
$ws.Range('E13').Value = '  +6.07%  '
$ws.Range('D14').Value = '1.798.60'
$ws.Range('E14').Value = '  +0.57%  '
$ws.Range('D15').Value = '0.624'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = '34.114.50'
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D18').Value = '68.02'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').Value = '243.69'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').Value = '0.0₃0778'
$ws.Range('E20').Value = '  -0.71%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').Value = '10.91'
$ws.Range('E22').Value = '  +2.24%  '
$ws.Range('E23').Value = '  +0.87%  '
$ws.Range('D25').Value = '161.85'
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').Value = '16.28'
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('E28').Value = '  +1.53%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  +2.06%  '
$ws.Range('D31').Value = '0.0518'
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('E33').Value = '  +4.11%  '
$ws.Range('E34').Value = '  +1.69%  '
$ws.Range('D35').Value = '1.407.22'
$ws.Range('E35').Value = '  +1.63%  '
$ws.Range('D36').Value = '0.655'
$ws.Range('E36').Value = '  +1.48%  '
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('E38').Value = '  +2.47%  '
$ws.Range('D39').Value = '2.35'
$ws.Range('E39').Value = '  +8.72%  '
$ws.Range('D40').Value = '80.19'
$ws.Range('E40').Value = '  +1.42%  '
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('E42').Value = '  +1.14%  '
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').Value = '13.35'
$ws.Range('E44').Value = '  +12.32%  '
$ws.Range('D45').Value = '0.0₆0141'
$ws.Range('E45').Value = '  +2.90%  '
$ws.Range('E46').Value = '  +4.49%  '
$ws.Range('B47').Value = 'Kaspa'
$ws.Range('C47').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D47').Value = '0.0507'
$ws.Range('E47').Value = '  +2.17%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '1.08'
$ws.Range('E48').Value = '  +2.51%  '
$ws.Range('D49').Value = '107.44'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').Value = '1.950.66'
$ws.Range('E50').Value = '  +0.09%  '